$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new ticker entry "GRT-USD" at the next empty row (row 73)
$ws.Range("A73").Value = "GRT-USD"
